$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-11-26 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-27 Wednesday", 2)

# Update the division problems inside the table. The table has 20 rows
# (5 data rows with content, interleaved with 3 blank rows each) and 5
# columns. Using the table's cells directly (rather than a global text
# Find/Replace) avoids ambiguity where the same problem text (e.g.
# "93÷7=") appears more than once but must be replaced with different
# values depending on its position.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="21÷4="},
    @{Row=1;  Col=2; New="63÷7="},
    @{Row=1;  Col=3; New="51÷7="},
    @{Row=1;  Col=4; New="17÷4="},
    @{Row=1;  Col=5; New="70÷7="},

    @{Row=5;  Col=1; New="44÷3="},
    @{Row=5;  Col=2; New="84÷6="},
    @{Row=5;  Col=3; New="23÷5="},
    @{Row=5;  Col=4; New="28÷9="},
    @{Row=5;  Col=5; New="72÷3="},

    @{Row=9;  Col=1; New="54÷7="},
    @{Row=9;  Col=2; New="39÷8="},
    @{Row=9;  Col=3; New="74÷2="},
    @{Row=9;  Col=4; New="53÷8="},
    @{Row=9;  Col=5; New="60÷8="},

    @{Row=13; Col=1; New="19÷7="},
    @{Row=13; Col=2; New="83÷3="},
    @{Row=13; Col=3; New="16÷5="},
    @{Row=13; Col=4; New="47÷8="},
    @{Row=13; Col=5; New="90÷6="},

    @{Row=17; Col=1; New="20÷4="},
    @{Row=17; Col=2; New="19÷9="},
    @{Row=17; Col=3; New="61÷9="},
    @{Row=17; Col=4; New="79÷3="},
    @{Row=17; Col=5; New="83÷5="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $u.New
}

Write-Host "Applied $($updates.Count) cell updates plus date update"
